# Insert a new data row before current row 83 (the "Macroferia Regional de
# Talca - Apio" weekly series), pushing the existing rows 83:210 down to
# 84:211 and keeping their formatting (incl. the date number format on
# column D).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(83).Insert()

# Populate the newly-inserted row 83 with the new weekly record.
$ws.Range("A83").Value = 5
$ws.Range("B83").Value = "Macroferia Regional de Talca"
$ws.Range("C83").Value = "Maule"
$ws.Range("D83").Value = 44775
$ws.Range("E83").Value = 7
$ws.Range("F83").Value = 100112017
$ws.Range("G83").Value = "Apio"
$ws.Range("H83").Value = "Americana (o)"
$ws.Range("I83").Value = "Primera"
$ws.Range("J83").Value = 600
$ws.Range("K83").Value = 8500
$ws.Range("L83").Value = 8500
$ws.Range("M83").Value = 8500
$ws.Range("N83").Value = "$/docena de matas"
$ws.Range("O83").Value = "Provincia del Elquí"
$ws.Range("P83").Value = 1417
$ws.Range("Q83").Value = 6
$ws.Range("R83").Value = "Hortaliza"
